$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "b"
$ws.Range("J3").Value = "Acknowledge (Backchannel)"
$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"
$ws.Range("I17").Value = "sd"
$ws.Range("J17").Value = "Statement-non-opinion"
$ws.Range("I18").Value = "sv"
$ws.Range("J18").Value = "Statement-opinion"
$ws.Range("I22").Value = "sv"
$ws.Range("J22").Value = "Statement-opinion"
$ws.Range("I23").Value = "b"
$ws.Range("J23").Value = "Acknowledge (Backchannel)"
$ws.Range("I26").Value = "sv"
$ws.Range("J26").Value = "Statement-opinion"
$ws.Range("I29").Value = "aa"
$ws.Range("J29").Value = "Agree/Accept"
$ws.Range("I32").Value = "sd"
$ws.Range("J32").Value = "Statement-non-opinion"
$ws.Range("I34").Value = "sd"
$ws.Range("J34").Value = "Statement-non-opinion"
$ws.Range("I38").Value = "qy"
$ws.Range("J38").Value = "Yes-No-Question"
$ws.Range("I39").Value = "b"
$ws.Range("J39").Value = "Acknowledge (Backchannel)"
$ws.Range("I42").Value = "sv"
$ws.Range("J42").Value = "Statement-opinion"
$ws.Range("I44").Value = "sv"
$ws.Range("J44").Value = "Statement-opinion"
$ws.Range("I46").Value = "aa"
$ws.Range("J46").Value = "Agree/Accept"
$ws.Range("I50").Value = "sv"
$ws.Range("J50").Value = "Statement-opinion"
$ws.Range("I55").Value = "b"
$ws.Range("J55").Value = "Acknowledge (Backchannel)"
$ws.Range("I57").Value = "qy"
$ws.Range("J57").Value = "Yes-No-Question"
$ws.Range("I65").Value = "ba"
$ws.Range("J65").Value = "Appreciation"
$ws.Range("I67").Value = "aa"
$ws.Range("J67").Value = "Agree/Accept"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I75").Value = "sd"
$ws.Range("J75").Value = "Statement-non-opinion"
$ws.Range("I83").Value = "aa"
$ws.Range("J83").Value = "Agree/Accept"
$ws.Range("I89").Value = "sv"
$ws.Range("J89").Value = "Statement-opinion"
$ws.Range("I104").Value = "sd"
$ws.Range("J104").Value = "Statement-non-opinion"
$ws.Range("I108").Value = "b"
$ws.Range("J108").Value = "Acknowledge (Backchannel)"
$ws.Range("I110").Value = "sd"
$ws.Range("J110").Value = "Statement-non-opinion"
$ws.Range("I118").Value = "sd"
$ws.Range("J118").Value = "Statement-non-opinion"
$ws.Range("I126").Value = "b"
$ws.Range("J126").Value = "Acknowledge (Backchannel)"
$ws.Range("I129").Value = "ba"
$ws.Range("J129").Value = "Appreciation"
$ws.Range("I135").Value = "aa"
$ws.Range("J135").Value = "Agree/Accept"
$ws.Range("I136").Value = "aa"
$ws.Range("J136").Value = "Agree/Accept"
$ws.Range("I138").Value = "%"
$ws.Range("J138").Value = "Uninterpretable"
$ws.Range("I142").Value = "ba"
$ws.Range("J142").Value = "Appreciation"
$ws.Range("I147").Value = "%"
$ws.Range("J147").Value = "Uninterpretable"
$ws.Range("I156").Value = "aa"
$ws.Range("J156").Value = "Agree/Accept"
$ws.Range("I157").Value = "ba"
$ws.Range("J157").Value = "Appreciation"
$ws.Range("I159").Value = "sd"
$ws.Range("J159").Value = "Statement-non-opinion"
$ws.Range("I172").Value = "sv"
$ws.Range("J172").Value = "Statement-opinion"
$ws.Range("I176").Value = "b"
$ws.Range("J176").Value = "Acknowledge (Backchannel)"
$ws.Range("I190").Value = "sv"
$ws.Range("J190").Value = "Statement-opinion"
$ws.Range("I197").Value = "b"
$ws.Range("J197").Value = "Acknowledge (Backchannel)"
$ws.Range("I223").Value = "%"
$ws.Range("J223").Value = "Uninterpretable"
$ws.Range("I226").Value = "qy"
$ws.Range("J226").Value = "Yes-No-Question"
$ws.Range("I242").Value = "sv"
$ws.Range("J242").Value = "Statement-opinion"
$ws.Range("I246").Value = "b"
$ws.Range("J246").Value = "Acknowledge (Backchannel)"
$ws.Range("I264").Value = "aa"
$ws.Range("J264").Value = "Agree/Accept"
$ws.Range("I270").Value = "sd"
$ws.Range("J270").Value = "Statement-non-opinion"
$ws.Range("I277").Value = "aa"
$ws.Range("J277").Value = "Agree/Accept"
$ws.Range("I284").Value = "sv"
$ws.Range("J284").Value = "Statement-opinion"
$ws.Range("I286").Value = "%"
$ws.Range("J286").Value = "Uninterpretable"
$ws.Range("I288").Value = "sv"
$ws.Range("J288").Value = "Statement-opinion"
$ws.Range("I299").Value = "%"
$ws.Range("J299").Value = "Uninterpretable"
$ws.Range("I305").Value = "aa"
$ws.Range("J305").Value = "Agree/Accept"
